$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 422.875
$ws.Range("J5").Value = 701
$ws.Range("L5").Value = 701
$ws.Range("N5").Value = -931

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 699.08
$ws.Range("I19").Value = 610.1
$ws.Range("J19").Value = 758.4
$ws.Range("K19").Value = 610.1
$ws.Range("L19").Value = 758.4
$ws.Range("M19").Value = -435.1
$ws.Range("N19").Value = -1108.4

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 262.9
$ws.Range("I33").Value = 269.8889
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 269.8889
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = -40.88889999999998
$ws.Range("N33").Value = -658

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 10101455
$ws.Range("I92").Value = 13889114
$ws.Range("J92").Value = 1033.3334
$ws.Range("K92").Value = 13889114
$ws.Range("L92").Value = 1033.3334
$ws.Range("M92").Value = -13887866
$ws.Range("N92").Value = -3529.3334

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 625297.3
$ws.Range("I98").Value = 625297.3
$ws.Range("K98").Value = 625297.3
$ws.Range("M98").Value = -623799.3

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2172.8462
$ws.Range("I116").Value = 2281
$ws.Range("K116").Value = 2281
$ws.Range("M116").Value = 1161

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 625297.3
$ws.Range("I122").Value = 625297.3
$ws.Range("K122").Value = 1875891.9
$ws.Range("M122").Value = -1873441.9

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1240.9
$ws.Range("I127").Value = 1066.6666
$ws.Range("J127").Value = 1315.5714
$ws.Range("K127").Value = 3199.9998
$ws.Range("L127").Value = 3946.7142
$ws.Range("M127").Value = 1760.0002
$ws.Range("N127").Value = -13866.7142

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 59784.39
$ws.Range("I132").Value = 63242.53
$ws.Range("K132").Value = 189727.59
$ws.Range("M132").Value = -187197.59

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2218.21
$ws.Range("I138").Value = 1273.6487
$ws.Range("J138").Value = 2772.9524
$ws.Range("K138").Value = 3820.9461
$ws.Range("L138").Value = 8318.8572
$ws.Range("M138").Value = 1319.0539
$ws.Range("N138").Value = -18598.8572

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 112913.445
$ws.Range("I2").Value = 126915.125
$ws.Range("J2").Value = 900
$ws.Range("K2").Value = 126915.125
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = -126802.125
$ws.Range("N2").Value = -1126

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16672.47
$ws.Range("I32").Value = 4959.5493
$ws.Range("J32").Value = 135103.11
$ws.Range("K32").Value = 4959.5493
$ws.Range("L32").Value = 135103.11
$ws.Range("M32").Value = -4672.5493
$ws.Range("N32").Value = -135677.11

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 112913.445
$ws.Range("I116").Value = 126915.125
$ws.Range("J116").Value = 900
$ws.Range("K116").Value = 126915.125
$ws.Range("L116").Value = 900
$ws.Range("M116").Value = -124621.125
$ws.Range("N116").Value = -5488

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1967.0476
$ws.Range("I122").Value = 1737.125
$ws.Range("J122").Value = 2702.8
$ws.Range("K122").Value = 5211.375
$ws.Range("L122").Value = 8108.400000000001
$ws.Range("M122").Value = -2761.375
$ws.Range("N122").Value = -13008.4

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 49285.668
$ws.Range("J139").Value = 49285.668
$ws.Range("L139").Value = 49285.668
$ws.Range("N139").Value = -59565.668

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 112913.445
$ws.Range("I3").Value = 126915.125
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 126915.125
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = -126801.125
$ws.Range("N3").Value = -1128

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1394.7407
$ws.Range("I16").Value = 1220.7059
$ws.Range("J16").Value = 1690.6
$ws.Range("K16").Value = 1220.7059
$ws.Range("L16").Value = 1690.6
$ws.Range("M16").Value = -933.7058999999999
$ws.Range("N16").Value = -2264.6

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3966.5068
$ws.Range("I31").Value = 1826.0605
$ws.Range("J31").Value = 5732.375
$ws.Range("K31").Value = 1826.0605
$ws.Range("L31").Value = 5732.375
$ws.Range("M31").Value = -1531.0605
$ws.Range("N31").Value = -6322.375

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3966.5068
$ws.Range("I34").Value = 1826.0605
$ws.Range("J34").Value = 5732.375
$ws.Range("K34").Value = 1826.0605
$ws.Range("L34").Value = 5732.375
$ws.Range("M34").Value = -1624.0605
$ws.Range("N34").Value = -6136.375

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 40002430
$ws.Range("I58").Value = 76924770
$ws.Range("J58").Value = 3236.6667
$ws.Range("K58").Value = 76924770
$ws.Range("L58").Value = 3236.6667
$ws.Range("M58").Value = -76924567
$ws.Range("N58").Value = -3642.6667

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1394.7407
$ws.Range("I113").Value = 1220.7059
$ws.Range("J113").Value = 1690.6
$ws.Range("K113").Value = 1220.7059
$ws.Range("L113").Value = 1690.6
$ws.Range("M113").Value = 949.2941000000001
$ws.Range("N113").Value = -6030.6

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6669510
$ws.Range("I132").Value = 9806015
$ws.Range("J132").Value = 4437.75
$ws.Range("K132").Value = 29418045
$ws.Range("L132").Value = 13313.25
$ws.Range("M132").Value = -29415515
$ws.Range("N132").Value = -18373.25

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 37502572
$ws.Range("I134").Value = 45456292
$ws.Range("J134").Value = 27781358
$ws.Range("K134").Value = 136368876
$ws.Range("L134").Value = 83344074
$ws.Range("M134").Value = -136366341
$ws.Range("N134").Value = -83349144

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 40002430
$ws.Range("I136").Value = 76924770
$ws.Range("J136").Value = 3236.6667
$ws.Range("K136").Value = 230774310
$ws.Range("L136").Value = 9710.000100000001
$ws.Range("M136").Value = -230771760
$ws.Range("N136").Value = -14810.0001

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2066.1035
$ws.Range("I5").Value = 1232.1818
$ws.Range("J5").Value = 4687
$ws.Range("K5").Value = 3696.5454
$ws.Range("L5").Value = 14061
$ws.Range("M5").Value = -3584.5454
$ws.Range("N5").Value = -14285

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 160
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 600
$ws.Range("N17").Value = -938
$ws.Range("M17").Value = -131

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2121.7778
$ws.Range("I129").Value = 5265
$ws.Range("J129").Value = 1223.7142
$ws.Range("K129").Value = 15795
$ws.Range("L129").Value = 3671.1426
$ws.Range("M129").Value = -10795
$ws.Range("N129").Value = -13671.1426

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7248072
$ws.Range("J131").Value = 7753685.5
$ws.Range("L131").Value = 23261056.5
$ws.Range("N131").Value = -23271136.5

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1814.2858
$ws.Range("I132").Value = 1750
$ws.Range("J132").Value = 1840
$ws.Range("K132").Value = 15750
$ws.Range("L132").Value = 16560
$ws.Range("M132").Value = -13220
$ws.Range("N132").Value = -21620

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2066.1035
$ws.Range("I135").Value = 1232.1818
$ws.Range("J135").Value = 4687
$ws.Range("K135").Value = 11089.6362
$ws.Range("L135").Value = 42183
$ws.Range("M135").Value = -8554.636200000001
$ws.Range("N135").Value = -47253

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2225.4375
$ws.Range("I122").Value = 2306.625
$ws.Range("J122").Value = 2144.25
$ws.Range("K122").Value = 6919.875
$ws.Range("L122").Value = 6432.75
$ws.Range("M122").Value = -4469.875
$ws.Range("N122").Value = -11332.75

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3008.6365
$ws.Range("I126").Value = 3040
$ws.Range("K126").Value = 9120
$ws.Range("M126").Value = -6650

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2979.8333
$ws.Range("I132").Value = 2668.0908
$ws.Range("K132").Value = 8004.2724
$ws.Range("M132").Value = -5474.2724

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3212.75
$ws.Range("I122").Value = 2343.6875
$ws.Range("J122").Value = 3908
$ws.Range("K122").Value = 7031.0625
$ws.Range("L122").Value = 11724
$ws.Range("M122").Value = -4581.0625
$ws.Range("N122").Value = -16624

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 56360.168
$ws.Range("I126").Value = 62927.438
$ws.Range("J126").Value = 3822
$ws.Range("K126").Value = 188782.314
$ws.Range("L126").Value = 11466
$ws.Range("M126").Value = -186312.314
$ws.Range("N126").Value = -16406

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3216.9387
$ws.Range("I132").Value = 3073.121
$ws.Range("J132").Value = 3513.5625
$ws.Range("K132").Value = 9219.363000000001
$ws.Range("L132").Value = 10540.6875
$ws.Range("M132").Value = -6689.363000000001
$ws.Range("N132").Value = -15600.6875
